$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results (case with 380 kV done)
$ws.Range("B2").Value = 19.46802728246466
$ws.Range("C2").Value = 10.33178790160724
$ws.Range("D2").Value = 6.801262432254918
$ws.Range("F2").Value = 36.29858851610221
$ws.Range("G2").Value = 3.683686426087322
$ws.Range("I2").Value = 28.89659979953459
$ws.Range("L2").Value = 10.534590121956
$ws.Range("M2").Value = 17.0089096745564
$ws.Range("N2").Value = 19.75516687194326
$ws.Range("B3").Value = 19.02707227529446
$ws.Range("C3").Value = 9.654239044558551
$ws.Range("D3").Value = 6.828634189735486
$ws.Range("F3").Value = 36.07150167452018
$ws.Range("G3").Value = 3.687538791747522
$ws.Range("I3").Value = 28.90435557035408
$ws.Range("L3").Value = 10.55001829051955
$ws.Range("M3").Value = 16.9379300108897
$ws.Range("N3").Value = 19.82548519306702
$ws.Range("B4").Value = 18.75679067686542
$ws.Range("C4").Value = 9.212527677963715
$ws.Range("D4").Value = 6.84623398843175
$ws.Range("F4").Value = 35.94281166093097
$ws.Range("G4").Value = 3.690026773236204
$ws.Range("I4").Value = 28.91627768399203
$ws.Range("L4").Value = 10.56094222412587
$ws.Range("M4").Value = 16.89796057029639
$ws.Range("N4").Value = 19.87062962820967
$ws.Range("B5").Value = 18.64693959545636
$ws.Range("C5").Value = 9.026048898148471
$ws.Range("D5").Value = 6.853606024353639
$ws.Range("F5").Value = 35.89310959121362
$ws.Range("G5").Value = 3.691071595794352
$ws.Range("I5").Value = 28.92293102162643
$ws.Range("L5").Value = 10.56575870501186
$ws.Range("M5").Value = 16.88259257622351
$ws.Range("N5").Value = 19.8895226060759
$ws.Range("B6").Value = 18.62872162339313
$ws.Range("C6").Value = 8.994692517986167
$ws.Range("D6").Value = 6.854842236961596
$ws.Range("F6").Value = 35.88502313510788
$ws.Range("G6").Value = 3.691246960539444
$ws.Range("I6").Value = 28.92414403668136
$ws.Range("L6").Value = 10.56658051656872
$ws.Range("M6").Value = 16.88009660463203
$ws.Range("N6").Value = 19.89268977880816
$ws.Range("B7").Value = 18.75530777688685
$ws.Range("C7").Value = 9.21003901009551
$ws.Range("D7").Value = 6.846332600007623
$ws.Range("F7").Value = 35.94213021926119
$ws.Range("G7").Value = 3.690040738609224
$ws.Range("I7").Value = 28.91636015353737
$ws.Range("L7").Value = 10.56100570343663
$ws.Range("M7").Value = 16.89774957330046
$ws.Range("N7").Value = 19.87088241466356
$ws.Range("B8").Value = 19.31599192322127
$ws.Range("C8").Value = 10.10349125469715
$ws.Range("D8").Value = 6.810535815624352
$ws.Range("F8").Value = 36.2180843399895
$ws.Range("G8").Value = 3.684989344391617
$ws.Range("I8").Value = 28.89778491488792
$ws.Range("L8").Value = 10.53960863145601
$ws.Range("M8").Value = 16.98369339749611
$ws.Range("N8").Value = 19.77900475549761
$ws.Range("B9").Value = 20.41181981585261
$ws.Range("C9").Value = 11.65249448763168
$ws.Range("D9").Value = 6.746612600277017
$ws.Range("F9").Value = 36.84251350509209
$ws.Range("G9").Value = 3.676051086714601
$ws.Range("I9").Value = 28.91837531959694
$ws.Range("L9").Value = 10.50916009006898
$ws.Range("M9").Value = 17.18033400874974
$ws.Range("N9").Value = 19.61439719436638
$ws.Range("B10").Value = 21.20527104797931
$ws.Range("C10").Value = 12.66793702272968
$ws.Range("D10").Value = 6.703442866474825
$ws.Range("F10").Value = 37.34912003796833
$ws.Range("G10").Value = 3.670066452867681
$ws.Range("I10").Value = 28.96848282906466
$ws.Range("L10").Value = 10.4938038060541
$ws.Range("M10").Value = 17.34112217868768
$ws.Range("N10").Value = 19.50286736639504
$ws.Range("B11").Value = 21.56186811619664
$ws.Range("C11").Value = 13.10352495183613
$ws.Range("D11").Value = 6.684621503959236
$ws.Range("F11").Value = 37.58926392494389
$ws.Range("G11").Value = 3.667468723694017
$ws.Range("I11").Value = 28.99889903256187
$ws.Range("L11").Value = 10.48833950478884
$ws.Range("M11").Value = 17.41761677983182
$ws.Range("N11").Value = 19.45415512471348
$ws.Range("B12").Value = 21.69613646259215
$ws.Range("C12").Value = 13.2647034380032
$ws.Range("D12").Value = 6.677611355393316
$ws.Range("F12").Value = 37.68152774696821
$ws.Range("G12").Value = 3.666502841019712
$ws.Range("I12").Value = 29.01151300622528
$ws.Range("L12").Value = 10.48648887457203
$ws.Range("M12").Value = 17.44704732796485
$ws.Range("N12").Value = 19.43599878981914
$ws.Range("B13").Value = 21.66725565171053
$ws.Range("C13").Value = 13.23015810376931
$ws.Range("D13").Value = 6.679115914110503
$ws.Range("F13").Value = 37.66159917022775
$ws.Range("G13").Value = 3.666710070325537
$ws.Range("I13").Value = 29.00874762801968
$ws.Range("L13").Value = 10.48687772286228
$ws.Range("M13").Value = 17.44068860781353
$ws.Range("N13").Value = 19.43989620226512
$ws.Range("B14").Value = 21.57293060893685
$ws.Range("C14").Value = 13.11686063559844
$ws.Range("D14").Value = 6.684042430100035
$ws.Range("F14").Value = 37.59682832935581
$ws.Range("G14").Value = 3.667388903427709
$ws.Range("I14").Value = 28.9999148301139
$ws.Range("L14").Value = 10.48818287284719
$ws.Range("M14").Value = 17.42002886254495
$ws.Range("N14").Value = 19.45265558709261
$ws.Range("B15").Value = 21.51504983302916
$ws.Range("C15").Value = 13.04697236137334
$ws.Range("D15").Value = 6.687075301976779
$ws.Range("F15").Value = 37.55732505628238
$ws.Range("G15").Value = 3.667807025797252
$ws.Range("I15").Value = 28.99464720013338
$ws.Range("L15").Value = 10.48901077442405
$ws.Range("M15").Value = 17.4074340059
$ws.Range("N15").Value = 19.46050880719386
$ws.Range("B16").Value = 21.18186674022035
$ws.Range("C16").Value = 12.63894185701067
$ws.Range("D16").Value = 6.704689291377819
$ws.Range("F16").Value = 37.33361537979275
$ws.Range("G16").Value = 3.670238720471294
$ws.Range("I16").Value = 28.96664845753676
$ws.Range("L16").Value = 10.49419150957799
$ws.Range("M16").Value = 17.33618897491421
$ws.Range("N16").Value = 19.50609145221951
$ws.Range("B17").Value = 20.97625200572912
$ws.Range("C17").Value = 12.38189482310079
$ws.Range("D17").Value = 6.715703838507177
$ws.Range("F17").Value = 37.19881137305419
$ws.Range("G17").Value = 3.671762347785747
$ws.Range("I17").Value = 28.95142466720072
$ws.Range("L17").Value = 10.49775925906116
$ws.Range("M17").Value = 17.29332810274559
$ws.Range("N17").Value = 19.53457240580034
$ws.Range("B18").Value = 20.85758649378352
$ws.Range("C18").Value = 12.23156760285278
$ws.Range("D18").Value = 6.72211601312569
$ws.Range("F18").Value = 37.12219154639579
$ws.Range("G18").Value = 3.672650442975233
$ws.Range("I18").Value = 28.94338584218347
$ws.Range("L18").Value = 10.49995454487652
$ws.Range("M18").Value = 17.26899224290236
$ws.Range("N18").Value = 19.55114439240927
$ws.Range("B19").Value = 20.81734394272313
$ws.Range("C19").Value = 12.18024278178293
$ws.Range("D19").Value = 6.724300283728669
$ws.Range("F19").Value = 37.09640866144062
$ws.Range("G19").Value = 3.672953157397675
$ws.Range("I19").Value = 28.94078724693198
$ws.Range("L19").Value = 10.50072243240735
$ws.Range("M19").Value = 17.26080745305894
$ws.Range("N19").Value = 19.55678813103265
$ws.Range("B20").Value = 20.99818264344842
$ws.Range("C20").Value = 12.40951451246173
$ws.Range("D20").Value = 6.714523365097365
$ws.Range("F20").Value = 37.21306717109058
$ws.Range("G20").Value = 3.671598940303158
$ws.Range("I20").Value = 28.95297100409898
$ws.Range("L20").Value = 10.49736464548199
$ws.Range("M20").Value = 17.29785808020796
$ws.Range("N20").Value = 19.53152085147863
$ws.Range("B21").Value = 21.6006580661188
$ws.Range("C21").Value = 13.15024099354156
$ws.Range("D21").Value = 6.682592217842751
$ws.Range("F21").Value = 37.61581763623251
$ws.Range("G21").Value = 3.667189030926533
$ws.Range("I21").Value = 29.00247949716455
$ws.Range("L21").Value = 10.48779358811761
$ws.Range("M21").Value = 17.42608469478805
$ws.Range("N21").Value = 19.44889998565366
$ws.Range("B22").Value = 21.98987965684061
$ws.Range("C22").Value = 13.61238752520178
$ws.Range("D22").Value = 6.662405718266761
$ws.Range("F22").Value = 37.88673774641437
$ws.Range("G22").Value = 3.664410722764631
$ws.Range("I22").Value = 29.0412237769767
$ws.Range("L22").Value = 10.4828122932618
$ws.Range("M22").Value = 17.51258202819113
$ws.Range("N22").Value = 19.39659197908232
$ws.Range("B23").Value = 21.78260373699763
$ws.Range("C23").Value = 13.36773443345929
$ws.Range("D23").Value = 6.673117312323773
$ws.Range("F23").Value = 37.74146072442677
$ws.Range("G23").Value = 3.665884095330032
$ws.Range("I23").Value = 29.01996105867083
$ws.Range("L23").Value = 10.48535441155136
$ws.Range("M23").Value = 17.46617650322389
$ws.Range("N23").Value = 19.42435548311155
$ws.Range("B24").Value = 20.98826920804928
$ws.Range("C24").Value = 12.39703558707667
$ws.Range("D24").Value = 6.715056808649461
$ws.Range("F24").Value = 37.20661937633463
$ws.Range("G24").Value = 3.671672778999622
$ws.Range("I24").Value = 28.95226968269746
$ws.Range("L24").Value = 10.49754260130654
$ws.Range("M24").Value = 17.29580912469143
$ws.Range("N24").Value = 19.53289984267533
$ws.Range("B25").Value = 20.11676399038944
$ws.Range("C25").Value = 11.25506653439012
$ws.Range("D25").Value = 6.763237001493562
$ws.Range("F25").Value = 36.66495827421989
$ws.Range("G25").Value = 3.678366322248318
$ws.Range("I25").Value = 28.90667404271576
$ws.Range("L25").Value = 10.51616505389817
$ws.Range("M25").Value = 17.12420890164664
$ws.Range("N25").Value = 19.65727009161293
